# Update Niger's (row 20) ISIN / Matched Udsteder / Matched Værdipapirets navn / Kommuner
# columns to an empty list, matching the refreshed exclusion-list data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D20").Value = "[]"
$ws.Range("E20").Value = "[]"
$ws.Range("F20").Value = "[]"
$ws.Range("G20").Value = "[]"
